$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.955.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.315.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0980"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.82%  "

$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.703.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.942.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.317.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("E19").Value = "  +2.50%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("E24").Value = "  -2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.992"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  +3.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("E32").Value = "  +4.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.924"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.89%  "

$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("E40").Value = "  +2.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.376"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.51%  "

$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "260.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0506"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.19%  "

$ws.Range("E47").Value = "  +1.78%  "

$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("E51").Value = "  +0.86%  "
